# Auto-generated Excel COM-interop script
# Updates the '想去人数' (F column) counts on each sheet
# to match the refreshed output data.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 201
$ws.Range("F4").Value = 382
$ws.Range("F7").Value = 6299
$ws.Range("F10").Value = 529
$ws.Range("F11").Value = 33
$ws.Range("F12").Value = 9549
$ws.Range("F14").Value = 2577
$ws.Range("F16").Value = 2367
$ws.Range("F17").Value = 2587
$ws.Range("F19").Value = 264
$ws.Range("F20").Value = 2031
$ws.Range("F23").Value = 354
$ws.Range("F27").Value = 52
$ws.Range("F30").Value = 1258
$ws.Range("F31").Value = 1233
$ws.Range("F32").Value = 86
$ws.Range("F33").Value = 114
$ws.Range("F35").Value = 1620
$ws.Range("F36").Value = 2698
$ws.Range("F38").Value = 959
$ws.Range("F39").Value = 333
$ws.Range("F40").Value = 1272
$ws.Range("F41").Value = 27
$ws.Range("F42").Value = 38

$ws = $wb.Worksheets.Item(2)
$ws.Range("F6").Value = 159
$ws.Range("F7").Value = 1
$ws.Range("F15").Value = 32
$ws.Range("F16").Value = 142

$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 717
$ws.Range("F3").Value = 933
$ws.Range("F4").Value = 108

$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 201
$ws.Range("F3").Value = 717
$ws.Range("F4").Value = 933
$ws.Range("F5").Value = 108
$ws.Range("F7").Value = 382
$ws.Range("F11").Value = 6299
$ws.Range("F13").Value = 529
$ws.Range("F14").Value = 33
$ws.Range("F15").Value = 9549
$ws.Range("F16").Value = 159
$ws.Range("F18").Value = 2577
$ws.Range("F20").Value = 2367
$ws.Range("F21").Value = 2587
$ws.Range("F23").Value = 264
$ws.Range("F24").Value = 2031
$ws.Range("F27").Value = 354
$ws.Range("F31").Value = 52
$ws.Range("F34").Value = 1258
$ws.Range("F35").Value = 1233
$ws.Range("F36").Value = 114
$ws.Range("F38").Value = 1620
$ws.Range("F40").Value = 2698
$ws.Range("F41").Value = 959
$ws.Range("F43").Value = 333
$ws.Range("F46").Value = 32
$ws.Range("F47").Value = 1272
$ws.Range("F48").Value = 38
$ws.Range("F50").Value = 142
$ws.Range("F51").Value = 142

